$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-18 Wednesday" "2024-12-19 Thursday"

Replace-Text "116×7=" "157×6="
Replace-Text "932×4=" "632×2="
Replace-Text "614×5=" "924×7="
Replace-Text "581×6=" "660×4="
Replace-Text "231×6=" "406×2="
Replace-Text "629×9=" "646×2="
Replace-Text "887×5=" "615×5="
Replace-Text "671×6=" "453×3="
Replace-Text "374×7=" "293×6="
Replace-Text "261×3=" "694×8="
Replace-Text "396×7=" "347×2="
Replace-Text "809×5=" "360×7="
Replace-Text "429×8=" "823×3="
Replace-Text "663×9=" "788×6="
Replace-Text "808×7=" "878×3="
Replace-Text "295×8=" "886×3="
Replace-Text "631×6=" "730×2="
Replace-Text "900×4=" "907×4="
Replace-Text "101×6=" "225×3="
Replace-Text "861×8=" "309×9="
Replace-Text "502×6=" "900×6="
Replace-Text "160×5=" "841×2="
Replace-Text "330×5=" "881×3="
Replace-Text "735×6=" "415×4="
Replace-Text "866×6=" "842×8="
